# Update de notas ADS-ETC-DCE
# Applies the 2nd-cuatrimestre grade updates on sheet "Cuatri B".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: add the missing K10 grade (propagates through L10/M10 formulas)
$ws.Range("K10").Value = 8.8

# Row 18: add the missing F18 grade (propagates through G18/K18 formulas)
$ws.Range("F18").Value = 3.71

# Row 22: clear C22 and H22 grades (propagates through D22/K22/L22 formulas)
$ws.Range("C22").Value = $null
$ws.Range("H22").Value = $null

# Row 30: L30 becomes a fixed value instead of the computed formula
$ws.Range("L30").Value = 0.35

# Move the active selection to N20 (scroll/selection state)
$ws.Range("N20").Select()
